$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the SW1 part row (row 27): the footprint description and the
# JLCPCB part number were revised.
$ws.Range("D27").Value = "C2921603"
$ws.Range("C27").Value = "SW_DIP_SPSTx01_Slide_Copal_CHS-01TA_W5.08mm_P1.27mm_Jpin"

# Update the active cell selection left in the sheet by the author.
$ws.Range("C20").Select()
